# Add team W/L/T record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: AD1="Wins", AE1="Losses", AF1="Ties"
# Copy formatting from the existing last header cell (AC1) so the new
# headers match the bold/centered/bordered header style (style index 1).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill every data row (2-43) with the team's record: 84 wins, 78 losses, 0 ties
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 30).Value = 84   # column AD - Wins
    $ws.Cells.Item($r, 31).Value = 78   # column AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # column AF - Ties
}
